{"js": "// Mark the \"2.2.2\" checklist item as done, the same way \"2.1.4\", \"2.1.5\"\n// and \"2.2.1\" were earlier: append \" X\" to the line, and (per the\n// reference edit) also give the paragraph a left tab stop at 1455 twips\n// (i.e. 72.75pt == 1455/20pt, since OOXML w:pos is in twentieths of a\n// point).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the paragraph whose whole text is exactly \"2.2.2\" (not e.g. the\n// later sentence that merely mentions \"2.2.1 , 2.4.5, 2.4.6\").\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"2.2.2\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error('Could not find the \"2.2.2\" paragraph');\n}\n\n// Pull the paragraph's own OOXML so we can splice in the <w:tabs> element\n// while keeping every existing attribute (rsids, etc.) untouched.\nconst ooxml = target.getOoxml();\nawait context.sync();\n\nlet xml = ooxml.value;\n\n// getOoxml() stamps synthetic Word-2010 paragraph/text ids on export;\n// the source paragraph never had them, so drop them again before we\n// feed the fragment back in via insertOoxml.\nxml = xml.replace(/\\s+w14:paraId=\"[0-9A-Fa-f]+\"/, \"\");\nxml = xml.replace(/\\s+w14:textId=\"[0-9A-Fa-f]+\"/, \"\");\n\n// Add the left tab stop at 1455 (twips) to this paragraph's pPr.\nxml = xml.replace(\n  \"<w:pPr><w:spacing\",\n  '<w:pPr><w:tabs><w:tab w:val=\"left\" w:pos=\"1455\"/></w:tabs><w:spacing'\n);\n\n// Append \" X\" to mark the item as finished, matching the commit\n// \"Xong \u0111\u1eb7t mua s\u1ea3n ph\u1ea9m 2.2.2\".\nxml = xml.replace(\"<w:t>2.2.2</w:t>\", \"<w:t>2.2.2 X</w:t>\");\n\ntarget.insertOoxml(xml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Mark the \"2.2.2\" checklist item as done, the same way \"2.1.4\", \"2.1.5\"\n# and \"2.2.1\" were earlier: append \" X\" to the line, and (per the\n# reference edit) also give the paragraph a left tab stop.\n#\n# Word's ParagraphFormat.TabStops.Add() takes the position in points,\n# while the OOXML w:tab/@w:pos is stored in twentieths of a point\n# (twips). The target w:pos=\"1455\" twips is therefore 1455/20 = 72.75 pt.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq \"2.2.2\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    throw 'Could not find the \"2.2.2\" paragraph'\n}\n\n$target.Range.ParagraphFormat.TabStops.Add(72.75)\n$target.Range.InsertAfter(\" X\")\n"}
